# admiral_cheatsheet: derive_extreme_records() code sample —
# remove the stray quotation marks around `dataset_add = "adlb"` so it
# reads `dataset_add = adlb` (matches the other unquoted dataset args).

$p = $ppt.ActivePresentation

$quote = [char]34
$oldSnippet = "  dataset = adlb,  dataset_add = " + $quote + "adlb" + $quote + ","
$newSnippet = "  dataset = adlb,  dataset_add = adlb,"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($hi = 1; $hi -le $slide.Shapes.Count; $hi++) {
        $shp = $slide.Shapes.Item($hi)
        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            $full = $tr.Text
            $idx = $full.IndexOf($oldSnippet)
            if ($idx -ge 0) {
                # Grab just the affected span as its own sub-range so the
                # surrounding runs (and the rest of the deck) stay untouched,
                # then overwrite its text in place.
                $sub = $tr.Characters($idx + 1, $oldSnippet.Length)
                $sub.Text = $newSnippet
            }
        }
    }
}
